# TC_138.xlsx edit script
# Applies the changes described by the commit diff:
#  1. Rename sheet "My Series" -> "Data"
#  2. Update cached CEIC add-in comment blob on A1 (opaque cache, set verbatim)
#  3. Update number format 166 "0.000" -> "###0.000"
#  4. D1: drop "SAR (China)" from the series name
#  5. A11: "Function Description" -> "Function Information"
#  6. D14: date serial 43778 -> 41781
#  7. C21/D21/E21: tiny floating point tweaks (Kurtosis row)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename worksheet
$ws.Name = "Data"

# 4. D1 series title text
$ws.Range("D1").Value = "(DC)Hong Kong Retail Bonds: Price: Mid: HK Link A: 07-05-2009: 3.60%"

# 5. A11 label text
$ws.Range("A11").Value = "Function Information"

# 6. D14 date serial value
$ws.Range("D14").Value = 41781

# 7. Kurtosis row tiny float corrections
$ws.Range("C21").Value = -1.268580586263763
$ws.Range("D21").Value = 2.931072119614401
$ws.Range("E21").Value = 1.8060927885682

# 3. Number format 166 -> "###0.000"
$ws.Range("B27").NumberFormat = "###0.000"

# 2. Cached add-in comment blob on A1 (opaque, set verbatim from target XML)
$newCommentText = "g0AAAB+LCAAAAAAAAAPtnOtvG1UWwP+VUaSVQKozM3bixNbtID/y8DZOQuw2LV+q8cxNMpvxjHceSf2t1QpaWApiKSBe210QsIsoBakLpQnd/wXVbvqp/8Kee+887XHqScuKIqMKcs/jPs6995zfTDOgly61dW4PW7ZmGqenxGlhisOGYqqasX16ynW2MmJ+6iUJLVxSsL4uW3IbO2DMgZdhFy/Z2umpHcfpFHl+f39/ej83bVrbfFYQRP58faWh7OC2nNEM25ENBU8FXuqTvaYkVFHbdezIquzIzPP0VK1Rm65gTamCrC4b8ja2psuurRnYthcMR3M0bBNPC8sOrlTr59jCpOx0flpE/JA8tCy7mq4yu5glk3t2MCxuam0sZQVxPiOKmWyhKcwVxdlirjA9k8u94jsGhmhFtp0GtvY0hQoajtzuUHdRzBaEOXE2V0B8ohH0FQZAQmu6uoH3NBurFazrdqqI8N4GlhQHVp0umALiI75eRyefwpIld3aamqPjtO6LpoUVCNSJxl7F+2uWF79mZwW0zR3NcrpVuZu6r7M2ttY6JBrpXCVUNQ2npGPLOduBTcUq7DkoJMdyMeJHKEOnqmYr8LNmuFiVtmTdjjrFlGjTtHbtjqzgVbiwPOlj39BNWYWT5Wi2oynhoEMKtG6ZHegRBi+buroIvXrGCYqg55oBISbDlk1zN5xdkhLRM0BPA+xpW3Z88yE5auyY+2uG3m24LVuxtBZWq2XfOlGHyM3zvCuu7ZhtmEUoQkwWkagqX6/zXfgH7tugElWxorVlfV2HUNpSDvqKCVDJdcwtzamYuts2bH9mA1K0Cetq4kvBOoM2WoMtNkjoTaNm+PYs2ImquMOGuR+MOaygoYiIS7bib/qwYtC4CjJ/E4c1dF/IKhc1HcpBdEci0vjZaOxg7CQeDKZBJPMtkgIjlburbrsFl6wFN22PjmojPtQjOK1w4mFekgBFI0P/NAWhSP/APAI1WjDU0Xa+EsFwkbEkcQbxAyIEa9LLumzsgnRTc3ZWS/5aEjSIRWCk/bAOwf3t6HKXioMoRWWoZii6q2KWFmrGFj2iZG5sU0eq0ZBoBW66hGSj2+x2IA3bWtGBH05PQWEu2o4FpX9KUkzXcKwuyR+I90yf5GO7LYMOIOtj+2xZ+M8uEEd30TWUiqmOP5rKonPW0JzxZ2i6FkuK47vQ6JH86NpVTDINTf1j+ytp1mRbqczbBm6bhqaMH20IMpm9eoKF2P6tGtsDs/s1tr0O1Z0VP3LXx3azABeh3KUapmTbpqLRw+pdDzXiz4+4MlW8Jbs6oJoDhXY7yL2DYlSydwdtoiJ01tL9DCgRELaBhBW1Pa0AQhDam1bMNhHwAKCbDcRH7QkHKXjB2F6RjW0XSCPIK4PyIP+SKtm0ZMMmywnAYiAVJxshP08x4JFY8lpz6UFgycsELeIH7FATtzumJet1CIy26B07j5oASOqys+O1oLbpWPGDzIeugVd8Zv7En2RGixRbBrnwXpocEFIjshbG3KFNKENklXW4lnpF1rWWxbKqX8qTdLBhISL6+ZcsLiUu+nsAz11Qfc/gLoHxsOHJ6ZEVfQU7wCSRSo2NmfnsTIE8YNA2oiteskzb5qpmGxPc44DpVFdxuBeWqusvQtmjQfEKzxj2UUu06OdyOtzLMKe4JG4ATLGtQakYNgw0oYP0sgvZB1t6N2LLVrRiKmDa+/jn3muvPrxx88G9r/qXD3qXDx8fXoM5Pj58na2KmaGm3NIxnU2zPD8v5KDUhyJEosdTAiaLpLILFyj4Bm3kPajRRmWhVllaKdOEEQh9d1YzePIM2DXdsNlg06cD0T3j/a1mJlLTT0BeO6aN1CCJPJPt4bh1VD/KkcXi4cFXDw9ujfT2AhbClFgoFDJC7omsBU+3YkbMjmCtlViSJ8YzGWE2k81GjAds0AbL8EGcaqqUE4WCKALVBclaDc5tktGgyuupKW/zA35MVGH8ExyBaNtX0lPfhHsRqNk9iDRYEPv/+evR7fdiVl50PUm8F5gc5RQyGO83aNerG02usXZ2o7LANRca5JyEuogd6/wYY2/04CZFD01w005xUMChYHFTddk6xf3RNfAproE7pzh4CpqK3MPYkUuSsuGequPBGUNacjssz0QcQmmCZZhwkjwSkg7V0dgmZZ9Qm+DhJaVbnyc5eGuphgArbdTLXMsIjimVoZiGiSJ6b4jrVx7c+++De/f6373Tu3sl1oM3TvAyAI49XK5oM7gFkAG9+jIgQZsNGs9d4WKknnhC8uC0bmqGY0uzWfrM5LUQuIqkN/pfVGtDiaMd03iBfECClmV74ZLj3XNpFfFxAcyzI0N1NcNnzUDAUnoY10ef/r3/8Z3++98dXf2698a/em++//Dg5tG3/2SXsP/ed/3r33pJf7Au0LmQJ1gGfRx9K6Jw5HJypFZzv1x+lzNMhwPE4FyaoH65/GGkMzJRCiNhz4BwwUTiUxgyjToTPy4ylWAOMb/AhRX8CqloucDCq2lmR1PCQV7JkK7I7aOKF2rNjGtjzgR6gno+YBw6j+vnufAKuz4nZMWsp2WzIUtoyXYk9Eu62ZJ1zlfQFw4DJjGv4x1CWzre0spaubQSmrBJrFkqtsgxZD8gHyFJhanZfss/ahEJaAH0FFcn74iGzIZVQc+RTMZ7L1u2SipJgsmvJmIWqOJaFoMjw3tV33A7AL/+i7nRevqyMsK7q4xNowQctmvVuB7aES3UxbiaCKiepiZPxdJUzSbvdRi+rpLQhE3QxV5wQji8t/EMvPY02BKe5J0FyzKtxOQTanyzOpAzZBQ+jHhgQ/eUUbYa7pUv8BPes6HtWVGYi9F2pVzkSuoe+esPu8itw0LZ25Fhzh5leQxhL49L2EOGSYTdwG1NNgxX1kdCdqX8+PCjozs/9r+53zv4of/Jzd61fz986yoIH9z9GmrQhLSfFWnPA2bnxyHt3P+FtIXZwpNJ2zOakPa4pL0s61uZLpaBKwMiBhg+KV+foLvRVL2cmqoTPY6l6hEZZxywvv5G76c7Y7L1+pkNrp3I1r5mBFs//PD+g7tXez/e6d2813//du/TL3vXP+vfvvGUkC2OhmzxYqSMDEN2fgLZE8h+viBbnED28wLZA8nnNwvZ8BQwH4PsF6qVF5dNuDRnyL82sCNrOlc2DZWAtAWwWOTqmlrkls9wK5qxy5WKnDBHUUgQCkUuN50X/pCE5M+m32MAvjo7JsAPGSYBfBWmN4rcHx9e6135pH/rM5a1Hx++/ujLD/t3bz/6+IfejW97l+/1rt0lGH/wY/8fh5Tnbz364A59i/6RZ/nu/f5bX5SgHa4RGpHgTbj/abkfYjqbEYDSZ57E/UKB7IGQ/9W5f34c7p+fcH8a7qcX9RQ3NcuZ4a+YnAD40/QzmvSrqUk/0eNY0h9MTWMgfv+DL8bk++Uz1SS4p+IRZA/5rPfqX56S47OjOT57MVKphjleFPLzE5KfkPxzRfLZCck/LyQ/kH5+syQ/K5C/NoyQ/OZ6jaI1pLEuV7Igm+jkbfiiaQJmL54rcufwdpFbsjA2uHUs20nUnr6PYwh9Mz8moQ8ZJhH6Jsa7IxG9//pPvbffYQTef/MqpD2AdlgMsDoF8oOjL6/03nut9+4VwuSf32Q/LJ4jr+RvfHP09idE/Onfjr5/bYLjzxDHRWGc1/DkHXxGnPm1f+EllxvjF16Y0QTHx8Vxdiu5Fob7apCKLTtcA86Gpcrd9FCevrfRaL6ZGs0TPY5F86GcNAabP7r8/ZhsDldohrwqmD0t0jo9ROkDBiN4fXQ36bE9NxrbcxcjZSkB2+cm798n1P58UXvud0Dt3ti/e2wfSD+/CrZ7P7DvGbwVmlWsYyf1B5m+d93cO7EvbH5a15q9pqteMNP9Qn0QlrCD6Fep5KCknUzJsgCjyAdsqT8jJX9lUEu3eknZyuYL87NzW/MFoZXLygVZaBW2WnheKAiFLNx11in5voJ00QbWJl8npBsE4DHuHn5JsiEb2yl7YwGmjuRrC4gUiTM5wpplO+dJnfR+YpILgeQCY+fz0gwD5POsfUHK5+eYBCz4aPd8bJ5+SnHYJ8umvqK1tZTfYAh+2ol3Ahvf6TC6TLuFUPJW8SVg30gPkKxbf4Jyxj5aStMbu0eQ4wN/8qGgrW3vOKnPlqJkW4KSzczkcD4zI8tzGRmLW/AAlJvLKzNzM7n8LPnM0OscEpqG91MOwvsbFv6vBKT/AU9rdtGDQAAA"
$ws.Range("A1").Comment.Text($newCommentText)
